$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.1554434735375247;  C = 0.3375848360084654;  D = 157.8057217802531;  E = 6.48142807727062;  G = 164.7801781670697 }
    3 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    4 = @{ B = 0.02258322285507441; C = 1.65323645889881;    D = 0.7127328510149897; E = 6.48142807727062;  G = 8.869980610039494 }
    5 = @{ B = 0.02258322285507441; C = 0.004309184025731883; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.239511964969853 }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    7 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
